$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2022-10-02)
$ws.Range("B2").Value = 0.003208871385164791
$ws.Range("C2").Value = 0.002571899574220771
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 4.037778955826875

# Row 3 (2022-08-06)
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 5.586269137925634

# Row 4 (2022-06-01)
$ws.Range("B4").Value = 0.6606524410359556
$ws.Range("C4").Value = 0.002571899574220771
$ws.Range("D4").Value = 0.1494219747398047
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 1.306882851410751
